# Insert a new daily-ranking record for 2026/01/19 at row 655, pushing the
# existing 2026/12/29 .. 2027/01/05 block down by one row (old 655-696 ->
# new 656-697), per the daily auto-push commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 655..696 down to 656..697, opening up a blank row 655.
$ws.Rows.Item(655).Insert()

# Column A stores dates as plain text (e.g. "2026/12/29"), not real date
# serials. Writing the literal string "2026/01/19" straight into .Value
# would get auto-parsed as a date. Instead, build it as a text-typed
# formula result and flatten it to a static value with PasteSpecial
# (xlPasteValues), which keeps the string as text without ever touching
# the cell's number format.
$ws.Range("A655").Formula = '="2026/01/19"'
$ws.Range("A655").Copy()
$ws.Range("A655").PasteSpecial(-4163)

$ws.Range("B655").Value = "月"
$ws.Range("C655").Value = 13
$ws.Range("D655").Value = 18
